# Weekly update: insert a new price record at the top of the
# "Macroferia Regional de Talca - Zanahoria" block (row 282), pushing the
# existing rows 282:383 down to 283:384.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 282 (shifts rows 282-383 down to 283-384,
# dimension grows from R383 to R384 automatically).
$ws.Rows(282).Insert()

# Populate the newly inserted row 282 with the new weekly record.
$ws.Range("A282").Value = 5
$ws.Range("B282").Value = "Macroferia Regional de Talca"
$ws.Range("C282").Value = "Maule"
$ws.Range("D282").Value = 44795
$ws.Range("E282").Value = 7
$ws.Range("F282").Value = 100114013
$ws.Range("G282").Value = "Zanahoria"
$ws.Range("H282").Value = "Sin especificar"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 300
$ws.Range("K282").Value = 10000
$ws.Range("L282").Value = 10000
$ws.Range("M282").Value = 10000
$ws.Range("N282").Value = "`$/saco 20 kilos"
$ws.Range("O282").Value = "Región de Ñuble"
$ws.Range("P282").Value = 500
$ws.Range("Q282").Value = 20
$ws.Range("R282").Value = "Hortaliza"
